$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B8").Value = "2025-08-20T17:48:34+01:00"
$ws.Range("B12").Value = "Code system for ovulation test results"
